$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 18.76993966666667
$ws.Range("H2").Value = 56.309819
$ws.Range("I2").Value = 0.1007685501185251
$ws.Range("J2").Value = 0.1007685501185251
$ws.Range("M2").Value = 15.03463666666667
$ws.Range("N2").Value = 45.10391
$ws.Range("O2").Value = 0.2402934356091235
$ws.Range("P2").Value = 0.2402934356091235
$ws.Range("Q2").Value = 282.1992231435877
$ws.Range("R2").Value = 2539.79300829229
$ws.Range("S2").Value = 0.02421402110933055
$ws.Range("T2").Value = 0.02421402110933055
$ws.Range("G3").Value = 18.76993966666667
$ws.Range("H3").Value = 56.309819
$ws.Range("I3").Value = 0.1007685501185251
$ws.Range("J3").Value = 0.1007685501185251
$ws.Range("O3").Value = 0.07715484716443403
$ws.Range("P3").Value = 0.07715484716443403
$ws.Range("Q3").Value = 90.61020695955644
$ws.Range("R3").Value = 815.4918626360079
$ws.Range("S3").Value = 0.007774782083376417
$ws.Range("T3").Value = 0.007774782083376417
$ws.Range("G4").Value = 18.76993966666667
$ws.Range("H4").Value = 56.309819
$ws.Range("I4").Value = 0.1007685501185251
$ws.Range("J4").Value = 0.1007685501185251
$ws.Range("M4").Value = 6.211932333333333
$ws.Range("N4").Value = 18.635797
$ws.Range("O4").Value = 0.09928318157880762
$ws.Range("P4").Value = 0.09928318157880764
$ws.Range("Q4").Value = 116.5975951100826
$ws.Range("R4").Value = 1049.378355990743
$ws.Range("S4").Value = 0.01000462225885071
$ws.Range("T4").Value = 0.01000462225885071
$ws.Range("G5").Value = 18.76993966666667
$ws.Range("H5").Value = 56.309819
$ws.Range("I5").Value = 0.1007685501185251
$ws.Range("J5").Value = 0.1007685501185251
$ws.Range("M5").Value = 36.49384133333334
$ws.Range("N5").Value = 109.481524
$ws.Range("O5").Value = 0.5832685356476348
$ws.Range("P5").Value = 0.5832685356476348
$ws.Range("Q5").Value = 684.9872000315729
$ws.Range("R5").Value = 6164.884800284156
$ws.Range("S5").Value = 0.05877512466696744
$ws.Range("T5").Value = 0.05877512466696744
$ws.Range("I6").Value = 0.5130361557055731
$ws.Range("J6").Value = 0.5130361557055731
$ws.Range("M6").Value = 15.03463666666667
$ws.Range("N6").Value = 45.10391
$ws.Range("O6").Value = 0.2402934356091235
$ws.Range("P6").Value = 0.2402934356091235
$ws.Range("Q6").Value = 1436.741963781313
$ws.Range("R6").Value = 12930.67767403182
$ws.Range("S6").Value = 0.1232792204461894
$ws.Range("T6").Value = 0.1232792204461894
$ws.Range("I7").Value = 0.5130361557055731
$ws.Range("J7").Value = 0.5130361557055731
$ws.Range("O7").Value = 0.07715484716443403
$ws.Range("P7").Value = 0.07715484716443403
$ws.Range("S7").Value = 0.03958322618329228
$ws.Range("T7").Value = 0.03958322618329228
$ws.Range("I8").Value = 0.5130361557055731
$ws.Range("J8").Value = 0.5130361557055731
$ws.Range("M8").Value = 6.211932333333333
$ws.Range("N8").Value = 18.635797
$ws.Range("O8").Value = 0.09928318157880762
$ws.Range("P8").Value = 0.09928318157880764
$ws.Range("Q8").Value = 593.6255100369327
$ws.Range("R8").Value = 5342.629590332393
$ws.Range("S8").Value = 0.05093586180340984
$ws.Range("T8").Value = 0.05093586180340984
$ws.Range("I9").Value = 0.5130361557055731
$ws.Range("J9").Value = 0.5130361557055731
$ws.Range("M9").Value = 36.49384133333334
$ws.Range("N9").Value = 109.481524
$ws.Range("O9").Value = 0.5832685356476348
$ws.Range("P9").Value = 0.5832685356476348
$ws.Range("Q9").Value = 3487.429355670739
$ws.Range("R9").Value = 31386.86420103665
$ws.Range("S9").Value = 0.2992378472726816
$ws.Range("T9").Value = 0.2992378472726816
$ws.Range("G10").Value = 20.061603
$ws.Range("H10").Value = 60.184809
$ws.Range("I10").Value = 0.1077029912330274
$ws.Range("J10").Value = 0.1077029912330274
$ws.Range("M10").Value = 15.03463666666667
$ws.Range("N10").Value = 45.10391
$ws.Range("O10").Value = 0.2402934356091235
$ws.Range("P10").Value = 0.2402934356091235
$ws.Range("Q10").Value = 301.61891205591
$ws.Range("R10").Value = 2714.57020850319
$ws.Range("S10").Value = 0.02588032178876347
$ws.Range("T10").Value = 0.02588032178876348
$ws.Range("G11").Value = 20.061603
$ws.Range("H11").Value = 60.184809
$ws.Range("I11").Value = 0.1077029912330274
$ws.Range("J11").Value = 0.1077029912330274
$ws.Range("O11").Value = 0.07715484716443403
$ws.Range("P11").Value = 0.07715484716443403
$ws.Range("Q11").Value = 96.845596312632
$ws.Range("R11").Value = 871.610366813688
$ws.Range("S11").Value = 0.008309807827736611
$ws.Range("T11").Value = 0.008309807827736611
$ws.Range("G12").Value = 20.061603
$ws.Range("H12").Value = 60.184809
$ws.Range("I12").Value = 0.1077029912330274
$ws.Range("J12").Value = 0.1077029912330274
$ws.Range("M12").Value = 6.211932333333333
$ws.Range("N12").Value = 18.635797
$ws.Range("O12").Value = 0.09928318157880762
$ws.Range("P12").Value = 0.09928318157880764
$ws.Range("Q12").Value = 124.621320334197
$ws.Range("R12").Value = 1121.591883007773
$ws.Range("S12").Value = 0.01069309563516939
$ws.Range("T12").Value = 0.01069309563516939
$ws.Range("G13").Value = 20.061603
$ws.Range("H13").Value = 60.184809
$ws.Range("I13").Value = 0.1077029912330274
$ws.Range("J13").Value = 0.1077029912330274
$ws.Range("M13").Value = 36.49384133333334
$ws.Range("N13").Value = 109.481524
$ws.Range("O13").Value = 0.5832685356476348
$ws.Range("P13").Value = 0.5832685356476348
$ws.Range("Q13").Value = 732.1249567743241
$ws.Range("R13").Value = 6589.124610968916
$ws.Range("S13").Value = 0.06281976598135797
$ws.Range("T13").Value = 0.06281976598135797
$ws.Range("G14").Value = 51.87415833333333
$ws.Range("H14").Value = 155.622475
$ws.Range("I14").Value = 0.2784923029428744
$ws.Range("J14").Value = 0.2784923029428744
$ws.Range("M14").Value = 15.03463666666667
$ws.Range("N14").Value = 45.10391
$ws.Range("O14").Value = 0.2402934356091235
$ws.Range("P14").Value = 0.2402934356091235
$ws.Range("Q14").Value = 779.9091229308056
$ws.Range("R14").Value = 7019.18210637725
$ws.Range("S14").Value = 0.0669198722648401
$ws.Range("T14").Value = 0.0669198722648401
$ws.Range("G15").Value = 51.87415833333333
$ws.Range("H15").Value = 155.622475
$ws.Range("I15").Value = 0.2784923029428744
$ws.Range("J15").Value = 0.2784923029428744
$ws.Range("O15").Value = 0.07715484716443403
$ws.Range("P15").Value = 0.07715484716443403
$ws.Range("Q15").Value = 250.4178652626889
$ws.Range("R15").Value = 2253.7607873642
$ws.Range("S15").Value = 0.02148703107002873
$ws.Range("T15").Value = 0.02148703107002873
$ws.Range("G16").Value = 51.87415833333333
$ws.Range("H16").Value = 155.622475
$ws.Range("I16").Value = 0.2784923029428744
$ws.Range("J16").Value = 0.2784923029428744
$ws.Range("M16").Value = 6.211932333333333
$ws.Range("N16").Value = 18.635797
$ws.Range("O16").Value = 0.09928318157880762
$ws.Range("P16").Value = 0.09928318157880764
$ws.Range("Q16").Value = 322.2387614152861
$ws.Range("R16").Value = 2900.148852737575
$ws.Range("S16").Value = 0.0276496018813777
$ws.Range("T16").Value = 0.0276496018813777
$ws.Range("G17").Value = 51.87415833333333
$ws.Range("H17").Value = 155.622475
$ws.Range("I17").Value = 0.2784923029428744
$ws.Range("J17").Value = 0.2784923029428744
$ws.Range("M17").Value = 36.49384133333334
$ws.Range("N17").Value = 109.481524
$ws.Range("O17").Value = 0.5832685356476348
$ws.Range("P17").Value = 0.5832685356476348
$ws.Range("Q17").Value = 1893.087303516878
$ws.Range("R17").Value = 17037.7857316519
$ws.Range("S17").Value = 0.1624357977266278
$ws.Range("T17").Value = 0.1624357977266278
